$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy existing date/time formatting down onto the two new rows so the
# new cells reuse the same style indexes (date style for col A, time
# style for col C) rather than minting new ones.
$ws.Range("A33").Copy()
$ws.Range("A34").PasteSpecial(-4122)
$ws.Range("A35").PasteSpecial(-4122)
$ws.Range("C33").Copy()
$ws.Range("C34").PasteSpecial(-4122)
$ws.Range("C35").PasteSpecial(-4122)

# New row 34: 2018-09-08, DualBoxing, 00:03:00
$ws.Cells.Item(34, 1).Value = 43351
$ws.Cells.Item(34, 2).Value = "DualBoxing"
$ws.Cells.Item(34, 3).Value = 0.0020833333333333333

# New row 35: 2018-09-16, ShadowOfWar, 00:37:20
$ws.Cells.Item(35, 1).Value = 43359
$ws.Cells.Item(35, 2).Value = "ShadowOfWar"
$ws.Cells.Item(35, 3).Value = 0.025925925925925925

# Update the saved view state: scroll position and the active selection.
$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Range("W41").Select()
